# Generate Report for Archive
# Update localization status from "Ready for handoff" to "In Translation"
# and shrink the corresponding status columns to their new auto-fit width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Closest value reachable through the ColumnWidth property (quantized in
# character-width units) to the narrower auto-fit width used for the
# status columns once their text shrank from "Ready for handoff" to
# "In Translation".
$newWidth = 12.5

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
